# Updated cryptos list on Fri Jul 28 17:14:16 UTC 2023 with GitHub Actions
# Refresh price (D) / volume-change (E) figures, and re-sort a few coins
# (rows 47-51) to match the new ranking order, per upstream data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    # Preserve these cells as plain text (they hold things like
    # "1.000" / "13.20") so Excel does not silently coerce the
    # string into a Double and eat the significant trailing zero.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "29.304.54"
Set-TextCell "E2" "  +0.07%  "

# Row 3
Set-TextCell "D3" "1.873.08"
Set-TextCell "E3" "  +0.31%  "

# Row 4
Set-TextCell "D4" "1.000"
Set-TextCell "E4" "  -0.04%  "

# Row 5
Set-TextCell "D5" "0.7115"
Set-TextCell "E5" "  -1.64%  "

# Row 6
Set-TextCell "D6" "241.77"
Set-TextCell "E6" "  +0.34%  "

# Row 7
Set-TextCell "E7" "  -0.11%  "

# Row 8
Set-TextCell "D8" "0.3110"
Set-TextCell "E8" "  +0.61%  "

# Row 9
Set-TextCell "D9" "0.07698"
Set-TextCell "E9" "  -1.80%  "

# Row 10
Set-TextCell "D10" "25.34"
Set-TextCell "E10" "  +0.44%  "

# Row 11
Set-TextCell "D11" "0.08365"
Set-TextCell "E11" "  +1.45%  "

# Row 12
Set-TextCell "D12" "1.871.15"
Set-TextCell "E12" "  +1.06%  "

# Row 13
Set-TextCell "D13" "5.226"
Set-TextCell "E13" "  -0.21%  "

# Row 14
Set-TextCell "D14" "0.7119"
Set-TextCell "E14" "  -1.50%  "

# Row 15
Set-TextCell "D15" "91.43"
Set-TextCell "E15" "  +0.65%  "

# Row 16
Set-TextCell "D16" "29.319.50"

# Row 17
Set-TextCell "D17" "0.000008253"
Set-TextCell "E17" "  +5.60%  "

# Row 18
Set-TextCell "D18" "5.951"
Set-TextCell "E18" "  +1.55%  "

# Row 19
Set-TextCell "D19" "242.29"
Set-TextCell "E19" "  -0.64%  "

# Row 20
Set-TextCell "D20" "2.134.47"
Set-TextCell "E20" "  +3.30%  "

# Row 21
Set-TextCell "D21" "13.20"
Set-TextCell "E21" "  -0.22%  "

# Row 22
Set-TextCell "D22" "0.9992"
Set-TextCell "E22" "  -0.07%  "

# Row 23
Set-TextCell "D23" "7.853"
Set-TextCell "E23" "  -1.77%  "

# Row 24
Set-TextCell "D24" "1.000"
Set-TextCell "E24" "  -0.06%  "

# Row 25
Set-TextCell "D25" "0.1629"
Set-TextCell "E25" "  +2.46%  "

# Row 26
Set-TextCell "D26" "163.30"
Set-TextCell "E26" "  +0.69%  "

# Row 27
Set-TextCell "D27" "9.011"
Set-TextCell "E27" "  +0.51%  "

# Row 28
Set-TextCell "D28" "18.51"
Set-TextCell "E28" "  +1.36%  "

# Row 29
Set-TextCell "E29" "  +0.59%  "

# Row 30
Set-TextCell "D30" "4.410"
Set-TextCell "E30" "  +0.08%  "

# Row 31
Set-TextCell "D31" "4.331"
Set-TextCell "E31" "  +5.49%  "

# Row 32
Set-TextCell "D32" "1.281"
Set-TextCell "E32" "  -4.84%  "

# Row 33
Set-TextCell "D33" "0.05250"
Set-TextCell "E33" "  +0.82%  "

# Row 34
Set-TextCell "D34" "1.930"

# Row 35
Set-TextCell "D35" "0.7541"
Set-TextCell "E35" "  +3.49%  "

# Row 36
Set-TextCell "D36" "1.174"
Set-TextCell "E36" "  -1.11%  "

# Row 37
Set-TextCell "D37" "2.680"
Set-TextCell "E37" "  -0.07%  "

# Row 38
Set-TextCell "D38" "0.01858"
Set-TextCell "E38" "  +0.03%  "

# Row 39
Set-TextCell "D39" "2.717"
Set-TextCell "E39" "  +0.63%  "

# Row 40
Set-TextCell "D40" "1.153.06"
Set-TextCell "E40" "  -1.50%  "

# Row 41
Set-TextCell "D41" "6.353"
Set-TextCell "E41" "  +4.29%  "

# Row 42
Set-TextCell "D42" "73.09"
Set-TextCell "E42" "  +1.27%  "

# Row 43
Set-TextCell "D43" "0.8883"
Set-TextCell "E43" "  -1.77%  "

# Row 44
Set-TextCell "D44" "104.76"

# Row 45
Set-TextCell "D45" "0.9996"
Set-TextCell "E45" "  -0.13%  "

# Row 46
Set-TextCell "D46" "2.029.78"
Set-TextCell "E46" "  +1.63%  "

# Row 47
Set-TextCell "B47" "Mantle"
Set-TextCell "C47" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D47" "0.5197"
Set-TextCell "E47" "  -1.67%  "

# Row 48
Set-TextCell "B48" "RenderToken"
Set-TextCell "C48" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D48" "1.791"
Set-TextCell "E48" "  +0.60%  "

# Row 49
Set-TextCell "D49" "9.379"
Set-TextCell "E49" "  +0.82%  "

# Row 50
Set-TextCell "B50" "TheSandbox"
Set-TextCell "C50" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell "D50" "0.4296"
Set-TextCell "E50" "  +0.60%  "

# Row 51
Set-TextCell "B51" "BabyDogeCoin"
Set-TextCell "C51" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D51" "0.00000000117"
Set-TextCell "E51" "  -2.47%  "
